# Adds the verso ("back side") of the Ketouva to the model: a page break
# followed by a right-aligned "${texteVerso}" merge-field placeholder,
# inserted right after the paragraph holding the "${fin}" placeholder
# (i.e. just before the section's closing sectPr).

$d = $word.ActiveDocument

# Locate the paragraph that holds the "${fin}" placeholder (normally the
# very last paragraph in the body) so the new content lands in the right
# spot even if the body already has trailing paragraphs.
$target = $null
ForEach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like '*${fin}*') {
        $target = $para
    }
}
if ($target -eq $null) {
    $target = $d.Paragraphs.Last
}

# Create a fresh empty paragraph right after it, then replace that empty
# paragraph's contents with the literal OOXML for the two new paragraphs
# (this preserves the "${fin}" paragraph untouched, unlike calling
# InsertXML directly on a collapsed end-of-story range, which eats the
# preceding paragraph mark).
$r = $target.Range
$r.Collapse(0)
$null = $r.InsertParagraphAfter()
$newRange = $d.Paragraphs.Last.Range

$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:bidi w:val="0"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Shlomo Stam" w:hAnsi="Shlomo Stam" w:cs="Shlomo Stam"/><w:spacing w:val="4"/><w:sz w:val="25"/><w:szCs w:val="25"/><w:lang w:val="fr-FR"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Shlomo Stam" w:hAnsi="Shlomo Stam" w:cs="Shlomo Stam"/><w:spacing w:val="4"/><w:sz w:val="25"/><w:szCs w:val="25"/><w:lang w:val="fr-FR"/></w:rPr><w:br w:type="page"/></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:suppressAutoHyphens/><w:spacing w:after="0" w:line="460" w:lineRule="exact"/><w:jc w:val="right"/><w:rPr><w:rFonts w:ascii="Shlomo Stam" w:hAnsi="Shlomo Stam" w:cs="Shlomo Stam"/><w:spacing w:val="4"/><w:sz w:val="30"/><w:szCs w:val="30"/><w:lang w:val="fr-FR"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Shlomo Stam" w:hAnsi="Shlomo Stam" w:cs="Shlomo Stam"/><w:spacing w:val="4"/><w:sz w:val="30"/><w:szCs w:val="30"/><w:lang w:val="fr-FR"/></w:rPr><w:lastRenderedPageBreak/><w:t>${</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Shlomo Stam" w:hAnsi="Shlomo Stam" w:cs="Shlomo Stam"/><w:spacing w:val="4"/><w:sz w:val="30"/><w:szCs w:val="30"/><w:lang w:val="fr-FR"/></w:rPr><w:t>texteVerso</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Shlomo Stam" w:hAnsi="Shlomo Stam" w:cs="Shlomo Stam"/><w:spacing w:val="4"/><w:sz w:val="30"/><w:szCs w:val="30"/><w:lang w:val="fr-FR"/></w:rPr><w:t>}</w:t></w:r></w:p>
'@

$null = $newRange.InsertXML($xml)
Write-Output "Inserted Ketouva verso (page break + texteVerso placeholder)."
